# "fix: update data dictionnary" -------------------------------------------
# Applies the corrections made to the data dictionary workbook:
#   - CLIENTUSER.password: fix FIELD SIZE (6 -> 50) and replace the masked
#     "******" EXAMPLE with a real example value "password123".
#   - BOOK.publication_date (and every other DATE field): DATA FORMAT label
#     changes from DD/MM/YYYY to YYYY/MM/DD (it is one shared value used by
#     every DATE row: publication_date, loan_date, excepted_return_date,
#     actual_return_date, reservation_date, date_notification).
#   - BOOK.publication_date EXAMPLE value updated to a newer date.
#   - LOAN.excepted_return_date: fix a typo in its CODE
#     (excepted_return_date -> excpected_return_date).
#   - LOAN/RESERVATION/NOTIFICATION date EXAMPLE values: converted from
#     plain text dates ("21/06/2024", ...) into real date values, displayed
#     with the new yyyy/mm/dd number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CLIENTUSER.password (row 11) ------------------------------------------
$ws.Range("E11").Value = 50
$ws.Range("G11").Value = "password123"

# --- Shared "DATA FORMAT" text for every DATE field (column F) -------------
# publication_date (19), loan_date (30), excepted_return_date (31),
# actual_return_date (32), reservation_date (36), date_notification (40).
$ws.Range("F19").Value = "YYYY/MM/DD"
$ws.Range("F30").Value = "YYYY/MM/DD"
$ws.Range("F31").Value = "YYYY/MM/DD"
$ws.Range("F32").Value = "YYYY/MM/DD"
$ws.Range("F36").Value = "YYYY/MM/DD"
$ws.Range("F40").Value = "YYYY/MM/DD"

# --- BOOK.publication_date EXAMPLE (row 19) ---------------------------------
$ws.Range("G19").Value = 45041

# --- LOAN entity (rows 30-32) -----------------------------------------------
# loan_date EXAMPLE: text date -> real date value with yyyy/mm/dd format.
$ws.Range("G30").Value = 45464
$ws.Range("G30").NumberFormat = "yyyy/mm/dd"

# excepted_return_date CODE typo fix.
$ws.Range("B31").Value = "excpected_return_date"

# excepted_return_date EXAMPLE: text date -> real date value.
$ws.Range("G31").Value = 45494
$ws.Range("G31").NumberFormat = "yyyy/mm/dd"

# actual_return_date EXAMPLE: text date -> real date value.
$ws.Range("G32").Value = 45498
$ws.Range("G32").NumberFormat = "yyyy/mm/dd"

# --- RESERVATION entity (row 36) --------------------------------------------
# reservation_date EXAMPLE: text date -> real date value.
$ws.Range("G36").Value = 45464
$ws.Range("G36").NumberFormat = "yyyy/mm/dd"

# --- NOTIFICATION entity (row 40) -------------------------------------------
# date_notification EXAMPLE: text date -> real date value.
$ws.Range("G40").Value = 45546
$ws.Range("G40").NumberFormat = "yyyy/mm/dd"
